$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.744.69'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.635.40'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.501'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.92%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.55'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.85%  '
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.24'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.38%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.861.20'
$ws.Range("E13").Value = '  -0.09%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.637.29'
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.554'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.29%  '
$ws.Range("D16").Value = '0.0₃0765'
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.73'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.76%  '
$ws.Range("D18").Value = '25.763.48'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("E20").Value = '  +1.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '193.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.29%  '
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.27'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.42%  '
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.03'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.52%  '
$ws.Range("E27").Value = '  -2.17%  '
$ws.Range("E29").Value = '  -0.62%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0494'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("E32").Value = '  +1.26%  '
$ws.Range("E33").Value = '  +0.48%  '
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("E35").Value = '  +0.38%  '
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.548'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.71%  '
$ws.Range("D38").Value = '1.117.55'
$ws.Range("E38").Value = '  -1.17%  '
$ws.Range("E39").Value = '  -1.88%  '
$ws.Range("E40").Value = '  -0.50%  '
$ws.Range("E41").Value = '  +0.76%  '
$ws.Range("E42").Value = '  +1.39%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.68'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.800'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.29%  '
$ws.Range("D45").Value = '1.769.28'
$ws.Range("E46").Value = '  -0.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.70%  '
$ws.Range("E48").Value = '  -2.27%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0502'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.24%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.62'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.76%  '
$ws.Range("E51").Value = '  +2.88%  '

Write-Host "Applied cryptos update"
